# fix extra time sum
# The "مجموع" (total) row sat at row 4, leaving an unused empty row above it.
# Delete that empty row so the totals row shifts up from row 4 to row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("گزارش تردد")

$ws.Rows("3:3").Delete()
